$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "25.949.00"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.644.23"
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "215.67"
$ws.Range("E5").Value = "  +0.27%  "
Set-TextValue "D6" "0.5056"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.27%  "
Set-TextValue "D8" "0.2581"
$ws.Range("E8").Value = "  +0.20%  "
Set-TextValue "D9" "0.06406"
$ws.Range("E9").Value = "  -0.72%  "
Set-TextValue "D10" "19.64"
$ws.Range("E10").Value = "  +0.86%  "
Set-TextValue "D11" "0.07794"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "4.291"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.623.89"
$ws.Range("E13").Value = "  -1.08%  "
Set-TextValue "D14" "0.5446"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "0.0₅7888"
$ws.Range("E15").Value = "  -0.38%  "
Set-TextValue "D16" "64.99"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "25.994.46"
$ws.Range("E18").Value = "  -0.27%  "
Set-TextValue "D19" "198.63"
$ws.Range("E19").Value = "  -2.18%  "
Set-TextValue "D20" "4.418"
Set-TextValue "D21" "9.986"
$ws.Range("E21").Value = "  -0.14%  "
Set-TextValue "D22" "6.008"
$ws.Range("E22").Value = "  +0.63%  "
Set-TextValue "D23" "1.006"
$ws.Range("E23").Value = "  -0.29%  "
Set-TextValue "D24" "1.872"
$ws.Range("E24").Value = "  -4.77%  "
Set-TextValue "D25" "141.00"
$ws.Range("E25").Value = "  -0.64%  "
Set-TextValue "D26" "0.1146"
$ws.Range("E26").Value = "  -0.43%  "
Set-TextValue "D27" "6.877"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +0.28%  "
Set-TextValue "D30" "0.05003"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  +0.71%  "
Set-TextValue "D32" "3.206"
$ws.Range("E32").Value = "  +0.28%  "
Set-TextValue "D33" "1.534"
$ws.Range("E33").Value = "  -0.38%  "
Set-TextValue "D34" "2.374"
$ws.Range("E34").Value = "  +1.20%  "
Set-TextValue "D35" "0.8960"
$ws.Range("E35").Value = "  +0.58%  "
Set-TextValue "D36" "2.614"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "1.146.01"
$ws.Range("E37").Value = "  +0.09%  "
Set-TextValue "D38" "0.5560"
$ws.Range("E38").Value = "  -1.05%  "
Set-TextValue "D39" "0.01564"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  -0.27%  "
Set-TextValue "D41" "5.700"
$ws.Range("E41").Value = "  +0.57%  "
Set-TextValue "D42" "0.8249"
$ws.Range("E42").Value = "  +2.11%  "
Set-TextValue "D43" "100.05"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").Value = "0.0₈122"
$ws.Range("E44").Value = "  +8.79%  "
$ws.Range("D45").Value = "1.782.78"
$ws.Range("E45").Value = "  +0.31%  "
Set-TextValue "D46" "0.4531"
$ws.Range("E46").Value = "  +0.09%  "
Set-TextValue "D47" "55.64"
$ws.Range("E47").Value = "  +1.19%  "
Set-TextValue "D48" "1.006"
$ws.Range("E48").Value = "  -0.25%  "
Set-TextValue "D49" "0.05063"
$ws.Range("E49").Value = "  +0.51%  "
Set-TextValue "D50" "1.006"
$ws.Range("E50").Value = "  -0.16%  "
Set-TextValue "D51" "0.09550"
$ws.Range("E51").Value = "  +2.65%  "
